# Natmi following Dr Hou advice
# Update the LR-pairs table (Ntf3-Ngfr) on Sheet1: the sending/target cluster
# set grows from {FAPs, sCs} to {ECs, FAPs, sCs}, which re-derives every
# summary statistic for every combination. Rewrite rows 2-7 in full.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) is untouched.

# Column layout (A..T):
# A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E Ligand-expressing cells, F Ligand detection rate,
# G Ligand average expression value, H Ligand total expression value,
# I Ligand derived specificity of average expression value,
# J Ligand derived specificity of total expression value,
# K Receptor-expressing cells, L Receptor detection rate,
# M Receptor average expression value, N Receptor total expression value,
# O Receptor derived specificity of average expression value,
# P Receptor derived specificity of total expression value,
# Q Edge average expression weight, R Edge total expression weight,
# S Edge average expression derived specificity,
# T Edge total expression derived specificity

$rows = @(
    @{ Row=2;  A="ECs";  D="FAPs"; E=2; F=0.6666666666666666; G=1.788130666666667;  H=5.364392;   I=0.08755714261138148; J=0.08755714261138148; K=3; L=1; M=3.362744666666666; N=10.088234; O=0.7488888671136141; P=0.748888867113614;  Q=6.013026862636444;   R=54.117241763728;   S=0.06557056933794263; T=0.06557056933794261 }
    @{ Row=3;  A="ECs";  D="sCs";  E=2; F=0.6666666666666666; G=1.788130666666667;  H=5.364392;   I=0.08755714261138148; J=0.08755714261138148; K=3; L=1; M=1.127567333333333; N=3.382702;  O=0.251111132886386;  P=0.2511111328863859; Q=2.016237727464889;   R=18.146139547184;   S=0.02198657327343886; T=0.02198657327343886 }
    @{ Row=4;  A="FAPs"; D="FAPs"; E=3; F=1;                  G=11.451921;          H=34.355763;  I=0.560751794520949;   J=0.560751794520949;   K=3; L=1; M=3.362744666666666; N=10.088234; O=0.7488888671136141; P=0.748888867113614;  Q=38.509886265838;     R=346.588976392542;  S=0.4199407761307197;  T=0.4199407761307196  }
    @{ Row=5;  A="FAPs"; D="sCs";  E=3; F=1;                  G=11.451921;          H=34.355763;  I=0.560751794520949;   J=0.560751794520949;   K=3; L=1; M=1.127567333333333; N=3.382702;  O=0.251111132886386;  P=0.2511111328863859; Q=12.912812023514;     R=116.215308211626;  S=0.1408110183902294;  T=0.1408110183902294  }
    @{ Row=6;  A="sCs";  D="FAPs"; E=3; F=1;                  G=7.182390333333333; H=21.547171;  I=0.3516910628676694;  J=0.3516910628676694;  K=3; L=1; M=3.362744666666666; N=10.088234; O=0.7488888671136141; P=0.748888867113614;  Q=24.15254478733489;   R=217.372903086014;  S=0.2633775216449518;  T=0.2633775216449518  }
    @{ Row=7;  A="sCs";  D="sCs";  E=3; F=1;                  G=7.182390333333333; H=21.547171;  I=0.3516910628676694;  J=0.3516910628676694;  K=3; L=1; M=1.127567333333333; N=3.382702;  O=0.251111132886386;  P=0.2511111328863859; Q=8.098628715115778;   R=72.887658436042;   S=0.08831354122271766; T=0.08831354122271765 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.A
    $ws.Range("B$row").Value = "Ntf3"
    $ws.Range("C$row").Value = "Ngfr"
    $ws.Range("D$row").Value = $r.D
    $ws.Range("E$row").Value = $r.E
    $ws.Range("F$row").Value = $r.F
    $ws.Range("G$row").Value = $r.G
    $ws.Range("H$row").Value = $r.H
    $ws.Range("I$row").Value = $r.I
    $ws.Range("J$row").Value = $r.J
    $ws.Range("K$row").Value = $r.K
    $ws.Range("L$row").Value = $r.L
    $ws.Range("M$row").Value = $r.M
    $ws.Range("N$row").Value = $r.N
    $ws.Range("O$row").Value = $r.O
    $ws.Range("P$row").Value = $r.P
    $ws.Range("Q$row").Value = $r.Q
    $ws.Range("R$row").Value = $r.R
    $ws.Range("S$row").Value = $r.S
    $ws.Range("T$row").Value = $r.T
}
